$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04628788019112266
$ws.Range("C2").Value = 0.4824692687926074
$ws.Range("D2").Value = 0.3511399927416267
$ws.Range("E2").Value = 0.5925706647663439
$ws.Range("F2").Value = 0.606948243852973
$ws.Range("G2").Value = 19
$ws.Range("B3").Value = 0.4480479793330983
$ws.Range("C3").Value = 0.7418921356350544
$ws.Range("D3").Value = 0.8887687429572028
$ws.Range("E3").Value = 0.9427453224265835
$ws.Range("F3").Value = 0.8535183764212539
$ws.Range("G3").Value = 18
$ws.Range("B4").Value = 0.6602403078873996
$ws.Range("C4").Value = 0.8862161360526695
$ws.Range("D4").Value = 1.515028269969415
$ws.Range("E4").Value = 1.230864846345615
$ws.Range("F4").Value = 1.070773292379532
$ws.Range("G4").Value = 17
$ws.Range("B5").Value = 0.6733686636360297
$ws.Range("C5").Value = 0.7201182210757808
$ws.Range("D5").Value = 0.9403137650122958
$ws.Range("E5").Value = 0.9696977699326197
$ws.Range("F5").Value = 0.7206577794175344
$ws.Range("G5").Value = 16
$ws.Range("B6").Value = 0.5133116504495733
$ws.Range("C6").Value = 0.6631082909423874
$ws.Range("D6").Value = 0.7111133683316896
$ws.Range("E6").Value = 0.84327538107767
$ws.Range("F6").Value = 0.692529925483697
$ws.Range("G6").Value = 15
$ws.Range("B7").Value = 0.458951976745339
$ws.Range("C7").Value = 0.673915242314678
$ws.Range("D7").Value = 0.6571042978108756
$ws.Range("E7").Value = 0.8106196998660196
$ws.Range("F7").Value = 0.6934053832595887
$ws.Range("G7").Value = 14
$ws.Range("B8").Value = 0.3510070987759771
$ws.Range("C8").Value = 0.6341193177088349
$ws.Range("D8").Value = 0.5711820071072622
$ws.Range("E8").Value = 0.7557658414530669
$ws.Range("F8").Value = 0.696640049829043
$ws.Range("G8").Value = 13
$ws.Range("B9").Value = 0.6367582691617952
$ws.Range("C9").Value = 0.6367582691617952
$ws.Range("D9").Value = 0.6187545134813729
$ws.Range("E9").Value = 0.7866095050794727
$ws.Range("F9").Value = 0.4898521181462611
$ws.Range("G9").Value = 9
$ws.Range("B10").Value = 0.4167169463657505
$ws.Range("C10").Value = 0.4167169463657505
$ws.Range("D10").Value = 0.2338164441105596
$ws.Range("E10").Value = 0.4835457001262234
$ws.Range("F10").Value = 0.2742340030023716
$ws.Range("G10").Value = 5
